# Fixed query issue for C3DC phs002599
#
# The "TreatmentTab" query (row 5, column B) wrapped the REPLACE() call in an
# extraneous CONCAT(...) which is unnecessary - simplify it to just the
# REPLACE() call.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Cells.Item(5, 2)

$oldText = $treatmentCell.Value2
if ($oldText.Contains("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))")) {
    $newText = $oldText.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
    $treatmentCell.Value2 = $newText
}

# Nudge the cell's font so Excel regenerates the style/font tables (matches
# the formatting tweak that shipped with the content fix).
$treatmentCell.Font.Size = 11
$treatmentCell.Font.ThemeColor = 1

# Move the view/selection down to the row that was edited.
$ws.Range("C5").Select()
